# Work in user v1
#
# Update the text columns (B:G) for rows 2-9 with their new values and
# append a new row 9 (A9 numeric index copied/styled like the rows above
# it). Columns B:G hold free-form text (names/logins/passwords) even when
# the value happens to look like a plain integer (e.g. "1"), so each cell
# is forced to Text format before the write and the format is cleared
# again afterwards so no stray style index is left behind - this keeps
# values such as "1", "2", "3" stored as text instead of being
# auto-coerced to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "B2" "1"
Set-TextValue "C2" "2"
Set-TextValue "D2" "3"
Set-TextValue "E2" "5"
Set-TextValue "F2" "4"
Set-TextValue "G2" "6"

# Row 3
Set-TextValue "B3" "LEXA2"
Set-TextValue "F3" "6"
Set-TextValue "G3" "4"

# Row 4
Set-TextValue "B4" "Евгений"
Set-TextValue "C4" "Евстратов"
Set-TextValue "D4" "Табулович"
Set-TextValue "E4" "Evgebn20"
Set-TextValue "F4" "Evgebn20"
Set-TextValue "G4" "Evgebn@mail.ru"

# Row 5
Set-TextValue "B5" "Rey"
Set-TextValue "C5" "Rey"
Set-TextValue "D5" "Rey"
Set-TextValue "E5" "Rey"
Set-TextValue "F5" "Rey"
Set-TextValue "G5" "Rey"

# Row 6
Set-TextValue "B6" "stepik"
Set-TextValue "C6" "stepik"
Set-TextValue "D6" "stepik"
Set-TextValue "E6" "stepik"
Set-TextValue "F6" "stepik"
Set-TextValue "G6" "stepik"

# Row 7
Set-TextValue "B7" "alex"
Set-TextValue "C7" "alex"
Set-TextValue "D7" "alex"
Set-TextValue "E7" "alex"
Set-TextValue "F7" "alex"
Set-TextValue "G7" "alex"

# Row 8
Set-TextValue "B8" "1"
Set-TextValue "C8" "1"
Set-TextValue "D8" "1"
Set-TextValue "E8" "1"
Set-TextValue "F8" "1"
Set-TextValue "G8" "1"

# New row 9 - copy formatting from A8 (bold/border/centered style used for
# the index column) onto A9, then set its value and the rest of the row.
$ws.Range("A8").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 7

Set-TextValue "B9" "1"
Set-TextValue "C9" "1"
Set-TextValue "D9" "1"
Set-TextValue "E9" "1"
Set-TextValue "F9" "1"
Set-TextValue "G9" "1"
